# Computes the ISO 7064 mod-97-10 IBAN check digits for a given country
# code + BBAN (here BBAN == the CCC number already present in column J).
function Get-IbanCheckDigits($country, $bban) {
    $rearranged = $bban + $country + "00"
    $letters = "ABCDEFGHIJKLMNOPQRSTUVWXYZ"
    $remainder = 0
    for ($i = 0; $i -lt $rearranged.Length; $i++) {
        $ch = $rearranged.Substring($i, 1)
        if ($ch -match "^[0-9]$") {
            $d = [int]$ch
            $remainder = ($remainder * 10 + $d) % 97
        } else {
            $idx = $letters.IndexOf($ch.ToUpper())
            $val = $idx + 10
            $d1 = [math]::Floor($val / 10)
            $d2 = $val % 10
            $remainder = ($remainder * 10 + $d1) % 97
            $remainder = ($remainder * 10 + $d2) % 97
        }
    }
    $check = 98 - $remainder
    return "{0:D2}" -f $check
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contribuyente")

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Running per-initials counter used to build the generated e-mail addresses,
# e.g. "Vittorio Diez Otero" -> "vdo00@vehiculos2025.com".
$emailCounts = @{}

for ($r = 8; $r -le $lastRow; $r++) {

    $b = $ws.Cells.Item($r, 2).Value()   # Apellido1
    $c = $ws.Cells.Item($r, 3).Value()   # Apellido2
    $d = $ws.Cells.Item($r, 4).Value()   # Nombre

    if ([string]::IsNullOrEmpty($b) -or [string]::IsNullOrEmpty($c) -or [string]::IsNullOrEmpty($d)) {
        continue
    }

    $initials = ($d.Substring(0, 1) + $b.Substring(0, 1) + $c.Substring(0, 1)).ToLower()

    if ($emailCounts.ContainsKey($initials)) {
        $n = $emailCounts[$initials]
    } else {
        $n = 0
    }
    $emailCounts[$initials] = $n + 1

    $email = "{0}{1:D2}@vehiculos2025.com" -f $initials, $n

    # Column G = Email
    $ws.Cells.Item($r, 7).Value = $email

    # Column J = CCC, Column K = IBAN
    $jVal = $ws.Cells.Item($r, 10).Value()
    if (-not [string]::IsNullOrEmpty($jVal)) {
        $kVal = $ws.Cells.Item($r, 11).Value()
        $needsNewIban = $true
        if (-not [string]::IsNullOrEmpty($kVal)) {
            if ($kVal.Length -ge 2 -and $kVal.Substring(0, 2) -eq "ES") {
                $needsNewIban = $false
            }
        }
        if ($needsNewIban) {
            $check = Get-IbanCheckDigits "ES" $jVal
            $newIban = "ES" + $check + $jVal
            $ws.Cells.Item($r, 11).Value = $newIban
        }
    }
}
